# update gebied onbekend en uitrustingsgraad/niveau
#
# Insert a new row above the old "99999/99999" row (row 4) so it becomes
# row 5, and populate the new row 4 with the "gebied onbekend" record:
#   A4 = 99993
#   B4 = #NULL! (error value, percentage number format like the old B4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing row 4 (99999/99999) down to row 5.
$ws.Rows.Item(4).Insert()

# New row 4: gebied-code 99993 with an error ("uitrustingsgraad/niveau"
# could not be computed -> #NULL!).
$ws.Range("A4").Value = 99993
$ws.Range("B4").Value = "#NULL!"
$ws.Range("B4").NumberFormat = "0.00%"
